$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").NumberFormat = "@"
$ws.Range("C3").NumberFormat = "@"

$ws.Range("A3").Value = "Swapnil Badve"
$ws.Range("B3").Value = "9975640367"
$ws.Range("C3").Value = "2025-03-18"
$ws.Range("D3").Value = 5500
$ws.Range("E3").Value = 100
$ws.Range("F3").Value = "Buffalo"
